# Scheduled data refresh: update market-price derived columns (H-N)
# across the Leve profit tables on each job sheet.
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 4575.5557
$ws.Range("J40").Value = 2760
$ws.Range("L40").Value = 2760
$ws.Range("N40").Value = -3110
$ws.Range("H51").Value = 7470.8335
$ws.Range("I51").Value = 15862.625
$ws.Range("J51").Value = 3274.9375
$ws.Range("K51").Value = 15862.625
$ws.Range("L51").Value = 3274.9375
$ws.Range("M51").Value = -15378.625
$ws.Range("N51").Value = -4242.9375
$ws.Range("H58").Value = 1962469.2
$ws.Range("H64").Value = 61845.293
$ws.Range("J64").Value = 3374.7778
$ws.Range("L64").Value = 3374.7778
$ws.Range("N64").Value = -3870.7778
$ws.Range("H67").Value = 61845.293
$ws.Range("J67").Value = 3374.7778
$ws.Range("L67").Value = 3374.7778
$ws.Range("N67").Value = -5090.7778
$ws.Range("H96").Value = 1657.625
$ws.Range("I96").Value = 594.9231
$ws.Range("J96").Value = 2913.5454
$ws.Range("K96").Value = 1784.7693
$ws.Range("L96").Value = 8740.636200000001
$ws.Range("M96").Value = -411.7692999999999
$ws.Range("N96").Value = -11486.6362
$ws.Range("H116").Value = 2766.6667
$ws.Range("H137").Value = 1498.1945
$ws.Range("I137").Value = 1239.6923
$ws.Range("K137").Value = 3719.0769
$ws.Range("M137").Value = -1169.0769
$ws.Range("H138").Value = 3275.641
$ws.Range("J138").Value = 3611.61
$ws.Range("L138").Value = 10834.83
$ws.Range("N138").Value = -21114.83

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H13").Value = 0
$ws.Range("J13").Value = 0
$ws.Range("L13").Value = 0
$ws.Range("N13").ClearContents()
$ws.Range("H44").Value = 5820
$ws.Range("J44").Value = 5820
$ws.Range("L44").Value = 5820
$ws.Range("N44").Value = -6796
$ws.Range("H45").Value = 34633.3
$ws.Range("J45").Value = 2117.125
$ws.Range("L45").Value = 2117.125
$ws.Range("N45").Value = -2871.125
$ws.Range("H55").Value = 14560
$ws.Range("J55").Value = 15486.667
$ws.Range("L55").Value = 15486.667
$ws.Range("N55").Value = -16116.667
$ws.Range("H74").Value = 1274.9642
$ws.Range("I74").Value = 1329.421
$ws.Range("J74").Value = 1160
$ws.Range("K74").Value = 1329.421
$ws.Range("L74").Value = 1160
$ws.Range("M74").Value = -455.421
$ws.Range("N74").Value = -2908
$ws.Range("H77").Value = 1274.9642
$ws.Range("I77").Value = 1329.421
$ws.Range("J77").Value = 1160
$ws.Range("K77").Value = 6647.105
$ws.Range("L77").Value = 5800
$ws.Range("M77").Value = -2279.105
$ws.Range("N77").Value = -14536
$ws.Range("H122").Value = 2409.04
$ws.Range("I122").Value = 2007.25
$ws.Range("J122").Value = 3123.3333
$ws.Range("K122").Value = 6021.75
$ws.Range("L122").Value = 9369.999899999999
$ws.Range("M122").Value = -3571.75
$ws.Range("N122").Value = -14269.9999

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 48823.047
$ws.Range("I20").Value = 56152.527
$ws.Range("J20").Value = 2403
$ws.Range("K20").Value = 56152.527
$ws.Range("L20").Value = 2403
$ws.Range("M20").Value = -55905.527
$ws.Range("N20").Value = -2897
$ws.Range("H105").Value = 183493.73
$ws.Range("I105").Value = 201850
$ws.Range("J105").Value = 168196.83
$ws.Range("K105").Value = 201850
$ws.Range("L105").Value = 168196.83
$ws.Range("M105").Value = -200103
$ws.Range("N105").Value = -171690.83
$ws.Range("H134").Value = 2932.7551
$ws.Range("I134").Value = 2788.2563
$ws.Range("J134").Value = 3496.3
$ws.Range("K134").Value = 8364.768899999999
$ws.Range("L134").Value = 10488.9
$ws.Range("M134").Value = -5829.768899999999
$ws.Range("N134").Value = -15558.9

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("M3:N3").ClearContents()
$ws.Range("H31").Value = 26736.725
$ws.Range("I31").Value = 858.4375
$ws.Range("J31").Value = 58586.92
$ws.Range("K31").Value = 858.4375
$ws.Range("L31").Value = 58586.92
$ws.Range("M31").Value = -563.4375
$ws.Range("N31").Value = -59176.92
$ws.Range("H34").Value = 26736.725
$ws.Range("I34").Value = 858.4375
$ws.Range("J34").Value = 58586.92
$ws.Range("K34").Value = 858.4375
$ws.Range("L34").Value = 58586.92
$ws.Range("M34").Value = -656.4375
$ws.Range("N34").Value = -58990.92
$ws.Range("H51").Value = 7966.3335
$ws.Range("J51").Value = 7966.3335
$ws.Range("L51").Value = 7966.3335
$ws.Range("N51").Value = -9438.333500000001
$ws.Range("H58").Value = 2393.842
$ws.Range("I58").Value = 2446.5454
$ws.Range("J58").Value = 2321.375
$ws.Range("K58").Value = 2446.5454
$ws.Range("L58").Value = 2321.375
$ws.Range("M58").Value = -2243.5454
$ws.Range("N58").Value = -2727.375
$ws.Range("H61").Value = 7966.3335
$ws.Range("J61").Value = 7966.3335
$ws.Range("L61").Value = 7966.3335
$ws.Range("N61").Value = -8662.333500000001
$ws.Range("H63").Value = 42990
$ws.Range("J63").Value = 42990
$ws.Range("L63").Value = 42990
$ws.Range("N63").Value = -44362
$ws.Range("H66").Value = 42990
$ws.Range("J66").Value = 42990
$ws.Range("L66").Value = 128970
$ws.Range("N66").Value = -135834
$ws.Range("H122").Value = 963.6
$ws.Range("I122").Value = 925
$ws.Range("J122").Value = 1021.5
$ws.Range("K122").Value = 2775
$ws.Range("L122").Value = 3064.5
$ws.Range("M122").Value = -325
$ws.Range("N122").Value = -7964.5
$ws.Range("H132").Value = 2470.353
$ws.Range("I132").Value = 2354.1936
$ws.Range("J132").Value = 3670.6667
$ws.Range("K132").Value = 7062.5808
$ws.Range("L132").Value = 11012.0001
$ws.Range("M132").Value = -4532.5808
$ws.Range("N132").Value = -16072.0001
$ws.Range("H136").Value = 2393.842
$ws.Range("I136").Value = 2446.5454
$ws.Range("J136").Value = 2321.375
$ws.Range("K136").Value = 7339.6362
$ws.Range("L136").Value = 6964.125
$ws.Range("M136").Value = -4789.6362
$ws.Range("N136").Value = -12064.125

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 680
$ws.Range("I122").Value = 530.8889
$ws.Range("J122").Value = 775.8570999999999
$ws.Range("K122").Value = 4778.0001
$ws.Range("L122").Value = 6982.7139
$ws.Range("M122").Value = -2328.0001
$ws.Range("N122").Value = -11882.7139
$ws.Range("H127").Value = 1186.1111
$ws.Range("J127").Value = 1230
$ws.Range("L127").Value = 3690
$ws.Range("N127").Value = -13610
$ws.Range("H131").Value = 1195.54
$ws.Range("I131").Value = 456.55
$ws.Range("J131").Value = 1380.2875
$ws.Range("K131").Value = 1369.65
$ws.Range("L131").Value = 4140.862499999999
$ws.Range("M131").Value = 3670.35
$ws.Range("N131").Value = -14220.8625

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 8300
$ws.Range("I57").Value = 2550
$ws.Range("K57").Value = 2550
$ws.Range("M57").Value = -1730
$ws.Range("H70").Value = 112810.9
$ws.Range("I70").Value = 160722.92
$ws.Range("K70").Value = 160722.92
$ws.Range("M70").Value = -160452.92
$ws.Range("H73").Value = 112810.9
$ws.Range("I73").Value = 160722.92
$ws.Range("K73").Value = 160722.92
$ws.Range("M73").Value = -159786.92
$ws.Range("H80").Value = 125277060
$ws.Range("J80").Value = 2993.3333
$ws.Range("L80").Value = 2993.3333
$ws.Range("N80").Value = -4989.3333
$ws.Range("H83").Value = 125277060
$ws.Range("J83").Value = 2993.3333
$ws.Range("L83").Value = 14966.6665
$ws.Range("N83").Value = -24950.6665
$ws.Range("H122").Value = 3798.5
$ws.Range("I122").Value = 2896.6667
$ws.Range("J122").Value = 6504
$ws.Range("K122").Value = 8690.000100000001
$ws.Range("L122").Value = 19512
$ws.Range("M122").Value = -6240.000100000001
$ws.Range("N122").Value = -24412
$ws.Range("H126").Value = 3924909.5
$ws.Range("I126").Value = 3637.8
$ws.Range("J126").Value = 11767453
$ws.Range("K126").Value = 10913.4
$ws.Range("L126").Value = 35302359
$ws.Range("M126").Value = -8443.400000000001
$ws.Range("N126").Value = -35307299

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H56").Value = 25919.834
$ws.Range("I56").Value = 9196.333000000001
$ws.Range("K56").Value = 9196.333000000001
$ws.Range("M56").Value = -8482.333000000001
$ws.Range("H80").Value = 3957643.5
$ws.Range("J80").Value = 3957643.5
$ws.Range("L80").Value = 3957643.5
$ws.Range("N80").Value = -3959639.5
$ws.Range("H83").Value = 3957643.5
$ws.Range("J83").Value = 3957643.5
$ws.Range("L83").Value = 11872930.5
$ws.Range("N83").Value = -11882914.5
$ws.Range("H107").Value = 125769.25
$ws.Range("I107").Value = 850.8
$ws.Range("J107").Value = 333966.66
$ws.Range("K107").Value = 2552.4
$ws.Range("L107").Value = 1001899.98
$ws.Range("M107").Value = -632.3999999999996
$ws.Range("N107").Value = -1005739.98

